$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "S"-size order row (row 3). This shifts row 4
# (phone number / total price) up into row 3.
$ws.Rows("3").Delete()

# The total price for the remaining order drops from 2000 to 1000.
$ws.Range("G3").Value = 1000
